# "Generate Report for Handoff"
# Updates the localization-status report: status flips from "In Translation"
# to "Ready for handoff", the associated timestamps advance, the stale
# "latest handback" commit hash in the error-detail message is refreshed,
# and the (now-wider) status columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps advance a little -----------------------------------------
$wsOverview.Range("G2").Value = "2017-02-09 13:48:07"
$wsDeDe.Range("H2").Value     = "2017-02-09 13:48:07"
$wsZhCn.Range("H2").Value     = "2017-02-09 13:47:47"

# --- Error detail message: refreshed "latest" handback commit hash -------
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d1b5de1721502e560b6602d8f4e47f0de9cc713a/e2e/24ccb9c7-e03d-4498-af8f-4682dacd4df7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13677b62b74b55ff260083e1a496a79e20134c39/e2e/24ccb9c7-e03d-4498-af8f-4682dacd4df7.md."

$wsZhCn.Range("R2").Value = $errorDetail
$wsDeDe.Range("R2").Value = $errorDetail

# --- Widen the status columns to fit "Ready for handoff" -----------------
# (16.333... characters is the closest input that lands on the same pixel
# grid as the target stored width of ~17.216 characters.)
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth     = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth     = 16.333333333333332
